$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# ---------------------------------------------------------------------------
# Settings sheet: remove the Orchestrator-queue related rows and replace them
# with the ACME credential / URL settings used by the non-Orchestrator flow.
# ---------------------------------------------------------------------------

# Old row 2 (OrchestratorQueueName / ProcessABCQueue / description) is removed.
$wsSettings.Range("A2:C2").Clear()

# Old row 3 (OrchestratorQueueFolder / ... ) is replaced with the
# logF_BusinessProcessName row that used to live on row 5.
$wsSettings.Range("A3:C3").Clear()
$wsSettings.Range("A3").Value = "logF_BusinessProcessName"
$wsSettings.Range("B3").Value = "Calculate_Client_Security_Hash"
$wsSettings.Range("C3").Value = "Logging field which allows grouping of log data of two or more subprocesses under the same business process name"
$wsSettings.Range("C3").WrapText = $true

# Row 5 now holds the ACME credential asset name.
$wsSettings.Range("A5:C5").Clear()
$wsSettings.Range("A5").Value = "ACME_Credentials"
$wsSettings.Range("B5").Value = "ACME_Credential"

# Row 6 now holds the ACME login URL, rendered as a hyperlink.
$wsSettings.Range("A6:C6").Clear()
$wsSettings.Range("A6").Value = "ACME_URL_Link"
$wsSettings.Range("B6").Value = "https://acme-test.uipath.com/login"
$wsSettings.Hyperlinks.Add($wsSettings.Range("B6"), "https://acme-test.uipath.com/login")

# ---------------------------------------------------------------------------
# Constants sheet: the process is no longer driven by Orchestrator queues, so
# bump the retry counters and flip ShouldMarkJobAsFaulted on.
# ---------------------------------------------------------------------------

$wsConstants.Range("B2").Value = 2
$wsConstants.Range("B3").Value = 3
$wsConstants.Range("B17").Value = $true

# ---------------------------------------------------------------------------
# Selection / active sheet state: the workbook was left with "Settings" as
# the active tab (instead of "Assets"), with new selections on each sheet.
# ---------------------------------------------------------------------------

$wsSettings.Activate()
$wsSettings.Range("A9").Select()

$wsConstants.Activate()
$wsConstants.Range("B18").Select()

$wsSettings.Activate()
